$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (interest count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 343
$wsExhibit.Range("F4").Value = 2912
$wsExhibit.Range("F6").Value = 612

# Sheet "全部类型" (All types) - update "想去人数" (interest count) column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 343
$wsAll.Range("F6").Value = 2912
$wsAll.Range("F8").Value = 612
